# Commit: change Ranker to ranker
# Lower-cases the "Ranker.com" domain references (both the DOMAIN column and
# the embedded JSON request bodies in the BODY column) for the "Ranker" test
# rows (13-23), and fixes row 21 ("...with wrong values"), which previously had
# stale copy/pasted DOMAIN/BODY content, to use the correct ranker.com values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13
$ws.Range("D13").Value = "ranker.com"
$ws.Range("E13").Value = "{`"x-ut-hb-params`":[ 
 {
 `"bidRequestId`": `"21b46f0d859b33`",
 `"domain`": `"ranker.com`",
 `"placementId`": `"10433394`",
 `"publisherId`": 3470,
 `"sizes`": [
 [1, 1],
 [160, 600]
 ],
 `"timeout`": 700,
 `"hbadaptor`": `"prebid`",
 `"params`": {}
}]}"

# Row 14
$ws.Range("D14").Value = "ranker.com"
$ws.Range("E14").Value = "{`"x-ut-hb-params`":[ 
 {
 `"bidRequestId`": `"21b46f0d859b33`",
 `"domain`": `"ranker.com`",
 `"placementId`": `"10433394`",
 `"sizes`": [
 [1, 1],
 [160, 600]
 ],
 `"timeout`": 700,
 `"hbadaptor`": `"prebid`",
 `"params`": {}
}]}"

# Row 15
$ws.Range("D15").Value = "ranker.com"
$ws.Range("E15").Value = "{`"x-ut-hb-params`":[ 
 {
 `"bidRequestId`": `"21b46f0d859b33`",
 `"domain`": `"ranker.com`",
 `"placementId`": `"10433394`",
 `"publisherId`": 3470,
 `"sizes`": [
 [1, 1],
 [160, 600]
 ],
 `"timeout`": 700,
 `"hbadaptor`": `"prebid`",
 `"params`": {`"placementId`" : `"10433394`",
                        `"publisherId`" : 3470
}
}]}"

# Row 16
$ws.Range("D16").Value = "ranker.com"
$ws.Range("E16").Value = "{`"x-ut-hb-params`":[ 
 {
 `"bidRequestId`": `"21b46f0d859b99`",
 `"domain`": `"ranker.com`",
 `"placementId`": `"10433394`",
 `"publisherId`": 3470,
 `"sizes`": [
 [1, 1],
 [160, 600]
 ],
 `"timeout`": 700,
 `"hbadaptor`": `"prebid`",
 `"params`": {`"placementId`" : `"10433394`",
                        `"publisherId`" : 3470
}
}]}"

# Row 17
$ws.Range("D17").Value = "ranker.com"
$ws.Range("E17").Value = "{`"x-ut-hb-params`":[ 
 {
 `"bidRequestId`": `"21b46f0d859b33`",
 `"placementId`": `"10433394`",
 `"publisherId`": 3470,
 `"sizes`": [
 [1, 1],
 [160, 600]
 ],
 `"timeout`": 700,
 `"hbadaptor`": `"prebid`",
 `"params`": {`"placementId`" : `"10433394`",
                        `"publisherId`" : 3470
}
}]}"

# Row 18
$ws.Range("D18").Value = "ranker.com"
$ws.Range("E18").Value = "{`"x-ut-hb-params`":[ 
 {
 `"bidRequestId`": `"21b46f0d859b33`",
 `"domain`": `"ranker.com`",
 `"publisherId`":3470,
 `"sizes`": [
 [1, 1],
 [160, 600]
 ],
 `"timeout`": 700,
 `"hbadaptor`": `"prebid`",
 `"params`": {`"publisherId`" : 3470
}
}]}"

# Row 19
$ws.Range("D19").Value = "ranker.com"
$ws.Range("E19").Value = "{`"x-ut-hb-params`":[ 
 {
 `"bidRequestId`": `"21b46f0d859b33`",
 `"domain`": `"ranker.com`",
 `"placementId`": `"10433394`",
 `"publisherId`": 3470,
 `"sizes`": [],
 `"timeout`": 700,
 `"hbadaptor`": `"prebid`",
 `"params`": {`"placementId`" : `"10433394`",
                        `"publisherId`" : 3470
}
}]}"

# Row 20
$ws.Range("D20").Value = "ranker.com"
$ws.Range("E20").Value = "{`"x-ut-hb-params`":[ 
 {
 `"bidRequestId`": `"21b46f0d859b33`",
 `"domain`": `"ranker.com`",
 `"placementId`": `"10433394`",
 `"publisherId`": 3470,
 `"sizes`": [
 [1, 1]
 ],
 `"timeout`": ,
 `"hbadaptor`": `"prebid`",
 `"params`": {`"placementId`" : `"10433394`",
                        `"publisherId`" : 3470
}
}]}"

# Row 21
$ws.Range("D21").Value = "ranker.com"
$ws.Range("E21").Value = "{`"x-ut-hb-params`":[ 
 {
 `"bidRequestId`": `"21b46f0d859b33`",
 `"domain`": `"ranker.com`",
 `"placementId`": `"10433394`",
 `"publisherId`": 0000,
 `"sizes`": [],
 `"timeout`": 700,
 `"hbadaptor`": `"prebida`",
 `"params`": {`"placementId`" : `"10433394`",
                        `"publisherId`" : 0000
}
}]}"

# Row 22
$ws.Range("D22").Value = "ranker.com"
$ws.Range("E22").Value = "{`"x-ut-hb-params`":[ 
 {
 `"bidRequestId`": `"21b46f0d859b47`",
 `"domain`": `"ranker.com`",
 `"placementId`": `"10433394`",
 `"publisherId`": 3470,
 `"sizes`": [
 [1, 1]
 ],
 `"timeout`": 700,
 `"hbadaptor`": `"prebid`",
 `"params`": {`"placementId`" : `"10433394`",
                        `"publisherId`" : 3470
}
}]}"

# Row 23
$ws.Range("D23").Value = "ranker.com"
$ws.Range("E23").Value = "{`"x-ut-hb-params`":[ 
 {
 `"bidRequestId`": `"21b46f0d859b47`",
 `"domain`": `"ranker.com`",
 `"placementId`": `"10433394`",
 `"publisherId`": 3470,
 `"sizes`": [
 [160, 600]
 ],
 `"timeout`": 700,
 `"hbadaptor`": `"prebid`",
 `"params`": {`"placementId`" : `"10433394`",
                        `"publisherId`" : 3470
}
}]}"

# Row 21 also had a stale taller row height from its old (longer) body text;
# the corrected body is shorter, so the row height shrinks to match row 19.
$ws.Rows.Item(21).RowHeight = 242.25
